$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-13 (columns A, B, C)
$data = @(
    @("5015948-83.2010.8.21.0001", "1999701-78.2010.8.21.0001", "Digitalizado"),
    @("5015948-83.2010.8.21.0001", "1999701-78.2010.8.21.0001", "Digitalizado"),
    @("5015948-83.2010.8.21.0001", "1999701-78.2010.8.21.0001", "Digitalizado"),
    @("5016829-79.2019.8.21.0022", "9000968-82.2019.8.21.0022", "Migrado"),
    @("5016824-57.2019.8.21.0022", "9007972-73.2019.8.21.0022", "Migrado"),
    @("5009671-07.2018.8.21.0022", "9009420-18.2018.8.21.0022", "Migrado"),
    @("5009668-52.2018.8.21.0022", "9009424-55.2018.8.21.0022", "Migrado"),
    @("5009667-67.2018.8.21.0022", "9008752-47.2018.8.21.0022", "Migrado"),
    @("5009665-97.2018.8.21.0022", "9009358-75.2018.8.21.0022", "Migrado"),
    @("5009663-30.2018.8.21.0022", "9007912-37.2018.8.21.0022", "Migrado"),
    @("5001397-28.2017.8.21.0042", "9001059-83.2017.8.21.0042", "Migrado"),
    @("5072049-91.2020.8.21.0001", "3221931-89.2005.8.21.0001", "Digitalizado")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}

# Remove the old rows 14-27 which are no longer part of the data set
$ws.Range("A14:C27").Clear()
